$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of dialogue lines, grouped by the same voice/name already
# used elsewhere in the sheet (folder grouping by voice name).
$ws.Range("A7").Value2 = 6
$ws.Range("B7").Value2 = $ws.Range("B6").Value2
$ws.Range("C7").Value2 = $ws.Range("C6").Value2
$ws.Range("D7").Value2 = "I am serious. And don't call me Shirley."

$ws.Range("A8").Value2 = 7
$ws.Range("B8").Value2 = $ws.Range("B5").Value2
$ws.Range("C8").Value2 = $ws.Range("C5").Value2
$ws.Range("D8").Value2 = "You talking to me?"

# Update the active selection to match the last-edited cell.
$ws.Range("D8").Select()

# Explicit page setup (portrait) as captured when the workbook was saved.
$ws.PageSetup.Orientation = 1

$wb.Save()
